# Updates cryptos list data (prices and volume percentages) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string updates (Coin names, Links, Volume%, and non-numeric-looking Prices) ---
$ws.Range("D2").Value = "72.201.61"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.910.21"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("E6").Value = "  +11.03%  "
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("E10").Value = "  +8.27%  "
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +5.64%  "
$ws.Range("D14").Value = "4.524.36"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "3.890.40"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  -5.11%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "71.907.20"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  -4.77%  "
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -6.21%  "
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("D35").Value = "0.0₃0994"
$ws.Range("E35").Value = "  +16.46%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -11.00%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +42.45%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("E45").Value = "  -6.77%  "
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("E48").Value = "  -16.13%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "2.878.38"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("E51").Value = "  +3.50%  "

# --- Price updates that look numeric: force as text to preserve original inlineStr formatting ---
# (set NumberFormat to Text before assigning, then restore default style afterward)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.672"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.772"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000326"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.129"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "607.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.424"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0473"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.145"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.32"
$ws.Range("D49").Style = "Normal"

Write-Output "Applied cryptos list update."
